$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. First paragraph: add two trailing spaces to the existing sentence,
#    then append a red parenthetical note typed as three separate runs
#    ("(This is a change - Ve" / "rsion for main branch" / ")").
# ---------------------------------------------------------------------

$firstPara = $d.Paragraphs(1).Range
$firstPara.End = $firstPara.End - 1
$firstPara.Select()
$word.Selection.Collapse(0)
$word.Selection.TypeText("  ")

$chunk1 = "(This is a change " + [char]0x2013 + " Ve"
$run1 = $d.Paragraphs(1).Range
$run1.End = $run1.End - 1
$run1.Select()
$word.Selection.Collapse(0)
$start1 = $word.Selection.Start
$word.Selection.TypeText($chunk1)
$d.Range($start1, $start1 + $chunk1.Length).Font.Color = 255

$chunk2 = "rsion for main branch"
$run2 = $d.Paragraphs(1).Range
$run2.End = $run2.End - 1
$run2.Select()
$word.Selection.Collapse(0)
$start2 = $word.Selection.Start
$word.Selection.TypeText($chunk2)
$d.Range($start2, $start2 + $chunk2.Length).Font.Color = 255

$chunk3 = ")"
$run3 = $d.Paragraphs(1).Range
$run3.End = $run3.End - 1
$run3.Select()
$word.Selection.Collapse(0)
$start3 = $word.Selection.Start
$word.Selection.TypeText($chunk3)
$d.Range($start3, $start3 + $chunk3.Length).Font.Color = 255

# ---------------------------------------------------------------------
# 2. Remove the trailing paragraph that reads
#    "ank God almighty, we are free at last." (the paragraph right
#    after "Shall be lifted-nevermore!").
# ---------------------------------------------------------------------

$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*God almighty, we are free at last*") {
        $target = $para
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
